# Add two new countries (Algeria, Andorra) to the Big Mac index table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 45: Algeria ------------------------------------------------------
# Columns A:C inherit the plain bordered style (style index 2) straight from
# the column's default format, same as every other data row in the sheet.
$ws.Range("A45").Value = "Algeria"
$ws.Range("B45").Value = 2.2
$ws.Range("C45").Value = 2.4

# --- Row 46: Andorra --------------------------------------------------------
# This row was pasted in with its own look: a 0.00 number format, a custom
# font/colour and a thin outline that only runs along the left/bottom (and,
# for the last column, the right) edges of each cell.
$ws.Range("A46").Value = "Andorra"

# B46 --------------------------------------------------------------------
$ws.Range("B46").Value = 5.29
$ws.Range("B46").Borders.Item(7).LineStyle = -4142   # xlEdgeLeft   -> none
$ws.Range("B46").Borders.Item(8).LineStyle = -4142   # xlEdgeTop    -> none
$ws.Range("B46").Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> none
$ws.Range("B46").Borders.Item(10).LineStyle = -4142  # xlEdgeRight  -> none
$ws.Range("B46").Borders.Item(7).LineStyle = 1        # left: thin
$ws.Range("B46").Borders.Item(9).LineStyle = 1        # bottom: thin
$ws.Range("B46").Interior.Color = 16777215
$ws.Range("B46").Font.Name = "Söhne"
$ws.Range("B46").Font.Size = 9.6
$ws.Range("B46").Font.Color = 855309
$ws.Range("B46").NumberFormat = "0.00"

# C46 --------------------------------------------------------------------
$ws.Range("C46").Value = 5.48
$ws.Range("C46").Borders.Item(7).LineStyle = -4142
$ws.Range("C46").Borders.Item(8).LineStyle = -4142
$ws.Range("C46").Borders.Item(9).LineStyle = -4142
$ws.Range("C46").Borders.Item(10).LineStyle = -4142
$ws.Range("C46").Borders.Item(7).LineStyle = 1        # left: thin
$ws.Range("C46").Borders.Item(9).LineStyle = 1        # bottom: thin
$ws.Range("C46").Borders.Item(10).LineStyle = 1       # right: thin
$ws.Range("C46").Interior.Color = 16777215
$ws.Range("C46").Font.Name = "Söhne"
$ws.Range("C46").Font.Size = 9.6
$ws.Range("C46").Font.Color = 855309
$ws.Range("C46").NumberFormat = "0.00"

# --- Restore the view: scrolled down so row 32 is at the top, with the
# cell just past the new data (D45) selected, matching where focus lands
# after typing the last value of the pasted block.
$excel.ActiveWindow.ScrollRow = 32
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D45").Select()
